$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_val data (filtered save games)
$values = @{
    2 = @(0.2917716402565462, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1, 1.242251378316819)
    3 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 5.586269137925634)
    4 = @(0.6606524410359556, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 0, 1.611132179096228)
    5 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 4.358119930609447)
    6 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 0, 8.974608811992548)
}

foreach ($row in $values.Keys) {
    $vals = $values[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
}
